$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '60.571.74'
Set-TextValue 'E2' '  -0.50%  '
Set-TextValue 'D3' '2.615.08'
Set-TextValue 'E3' '  -0.48%  '
Set-TextValue 'E4' '  +0.18%  '
Set-TextValue 'D5' '511.53'
Set-TextValue 'E5' '  +0.28%  '
Set-TextValue 'D6' '154.74'
Set-TextValue 'E6' '  -2.06%  '
Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  +0.36%  '
Set-TextValue 'E8' '  -2.76%  '
Set-TextValue 'D9' '2.628.70'
Set-TextValue 'E9' '  -1.39%  '
Set-TextValue 'D10' '6.70'
Set-TextValue 'E10' '  +4.72%  '
Set-TextValue 'E11' '  -0.50%  '
Set-TextValue 'E12' '  -0.33%  '
Set-TextValue 'E13' '  +1.46%  '
Set-TextValue 'D14' '3.075.86'
Set-TextValue 'E14' '  -0.75%  '
Set-TextValue 'D15' '60.494.37'
Set-TextValue 'E15' '  -0.13%  '
Set-TextValue 'E16' '  -1.08%  '
Set-TextValue 'D17' '0.0000141'
Set-TextValue 'E17' '  -0.14%  '
Set-TextValue 'D18' '2.624.07'
Set-TextValue 'E18' '  -1.35%  '
Set-TextValue 'E19' '  -0.72%  '
Set-TextValue 'D20' '351.38'
Set-TextValue 'E20' '  +1.38%  '
Set-TextValue 'D21' '10.61'
Set-TextValue 'E21' '  +0.65%  '
Set-TextValue 'D22' '6.18'
Set-TextValue 'E22' '  -0.36%  '
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  +0.30%  '
Set-TextValue 'D24' '60.61'
Set-TextValue 'E24' '  +0.59%  '
Set-TextValue 'E25' '  -0.16%  '
Set-TextValue 'E26' '  -0.71%  '
Set-TextValue 'E27' '  -0.04%  '
Set-TextValue 'D28' '0.0₃0844'
Set-TextValue 'E28' '  -3.26%  '
Set-TextValue 'D29' '7.39'
Set-TextValue 'E29' '  -2.31%  '
Set-TextValue 'E30' '  +0.24%  '
Set-TextValue 'D31' '19.46'
Set-TextValue 'E31' '  -0.48%  '
Set-TextValue 'E32' '  +0.09%  '
Set-TextValue 'D33' '150.60'
Set-TextValue 'E33' '  -4.19%  '
Set-TextValue 'D34' '5.82'
Set-TextValue 'E34' '  +0.81%  '
Set-TextValue 'E35' '  -2.04%  '
Set-TextValue 'E36' '  -2.23%  '
Set-TextValue 'D37' '0.891'
Set-TextValue 'E37' '  +5.33%  '
Set-TextValue 'E38' '  -0.78%  '
Set-TextValue 'D39' '0.849'
Set-TextValue 'E39' '  -1.14%  '
Set-TextValue 'D40' '36.40'
Set-TextValue 'E40' '  +3.39%  '
Set-TextValue 'D41' '3.78'
Set-TextValue 'E41' '  +0.19%  '
Set-TextValue 'D42' '294.99'
Set-TextValue 'E42' '  -5.43%  '
Set-TextValue 'E43' '  -2.99%  '
Set-TextValue 'E44' '  +0.07%  '
Set-TextValue 'E45' '  +0.19%  '
Set-TextValue 'B46' 'EnergySwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '19.92'
Set-TextValue 'E46' '  -1.29%  '
Set-TextValue 'B47' 'Hedera'
Set-TextValue 'C47' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D47' '0.0555'
Set-TextValue 'E47' '  -4.57%  '
Set-TextValue 'D48' '4.91'
Set-TextValue 'E48' '  -0.93%  '
Set-TextValue 'E49' '  -1.01%  '
Set-TextValue 'E50' '  +0.28%  '
Set-TextValue 'D51' '2.004.73'
Set-TextValue 'E51' '  -3.24%  '
